$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the extra blank row (row 2) so the data shifts up by one row.
$ws.Rows("2").Delete()

# Update the selected cell to match the new state.
$ws.Range("C6").Select()
